# Edit: extend the closing paragraph of "Experimentos propostos" with a
# transition sentence, then add two new paragraphs describing the test
# strategy (bogusMode) and the total-order multicast test, as per the
# commit "Descricao do teste total, remocao dos testes antigos."

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the target paragraph robustly by searching for the tail of its
# text, then expanding the found range to the whole paragraph.
$anchor = $d.Content
$null = $anchor.Find.Execute("comparando a saída tabulada de cada um.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Expand(4)  # wdParagraph

# --- Paragraph 1: the original paragraph, with the new sentence appended
# to its existing run (keeps the leading <w:tab/> as a distinct element).
$para1Text = "Como não seria necessário mais simular a performance de máquinas isoladas, apenas confirmar a sequencia de mensagens enviadas ou a exclusão mútua de execução, os experimentos desta vez foram todos realizados na mesma máquina, executando os processos manualmente e em seguida comparando a saída tabulada de cada um. Pela arquitetura escolhida os testes envolvem iniciar um conjunto de daemons usando o script t3daemons.lua e em seguida chamar os scripts que iniciam os processos, entram nos grupos e seguem o roteiro de comunicação para testar cada tipo de comunicação, conforme indicando a seguir."

$para1 = '<w:p ' + $wNs + '>' + `
    '<w:pPr><w:pStyle w:val="style18"/></w:pPr>' + `
    '<w:r><w:rPr/><w:tab/><w:t>' + $para1Text + '</w:t></w:r>' + `
  '</w:p>'

# --- Paragraph 2: strategy description, introducing the bogusMode parameter.
$para2 = '<w:p ' + $wNs + '>' + `
    '<w:pPr><w:pStyle w:val="style18"/></w:pPr>' + `
    '<w:r><w:rPr/><w:tab/><w:t xml:space="preserve">A principal estratégia de teste foi a sugerida na definição do trabalho. A biblioteca recebe em seu método de início um parâmetro, </w:t></w:r>' + `
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>bogusMode</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>, indicando se a biblioteca deve usar os métodos de efetivos (</w:t></w:r>' + `
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>bogusMode = false</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>, valor padrão) ou desviar para métodos que não fazem os controles (</w:t></w:r>' + `
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>bogusMode = true</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>). O roteiro de transmissão então irá fazer várias comunicações entre os processos, e o resultado quando a biblioteca estiver no modo efetivo deve ser igual ao esperado na definição de cada tipo de sincronização. Já no modo não efetivo deve ser possível encontrar situações onde as condições garantidas pelas sincronizações são violadas.</w:t></w:r>' + `
  '</w:p>'

# --- Paragraph 3: describes the total-order multicast test (alice/bob/carl).
$para3 = '<w:p ' + $wNs + '>' + `
    '<w:pPr><w:pStyle w:val="style0"/></w:pPr>' + `
    '<w:r><w:rPr/><w:tab/><w:t>Para testar o multicast total o script dos processos transmite uma mensagem contendo um nome designando cada processo (</w:t></w:r>' + `
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>alice, bob, carl</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>) e um número serial das mensagens. Ao receber uma mensagem da biblioteca o processo imprime a mensagem para um arquivo próprio (</w:t></w:r>' + `
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>alice.log, bob.log, carl.log</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>) e se a mensagem veio de outro processo o script aguarda entre 1 e 3 segundos (usando uma chamada não bloqueante aproveitando a função socket.select) e então envia outra mensagem para o grupo. A saída dos processos deve ter todas as mensagens recebidas por todos os processos na mesma ordem, mas quando a biblioteca for colocada no modo não efetivo devem existir mensagens recebidas em ordens diferentes.</w:t></w:r>' + `
  '</w:p>'

$anchor.InsertXML($para1 + $para2 + $para3)
